$d = $word.ActiveDocument

$find = $d.Content.Find

# Paragraph 1: "AccountDescription: {{AccountDescription}}" -> "{#rows}"
$find.Execute("AccountDescription: {{AccountDescription}}", $true, $false, $false, $false, $false, $true, 1, $false, "{#rows}", 2) | Out-Null

# Paragraph 3: "Request: {{Antrag}}" -> "Antrag: {Antrag}"
$find.Execute("Request: {{Antrag}}", $true, $false, $false, $false, $false, $true, 1, $false, "Antrag: {Antrag}", 2) | Out-Null

Write-Host "After text replacements:"
foreach ($p in $d.Paragraphs) {
    Write-Host "PARA: [" $p.Range.Text "]"
}

# Now fill the empty paragraph 2 with "Account Description: {AccountDescription}"
$p2 = $d.Paragraphs(2).Range
$p2.InsertBefore("Account Description: {AccountDescription}")

Write-Host "After filling paragraph 2:"
foreach ($p in $d.Paragraphs) {
    Write-Host "PARA: [" $p.Range.Text "]"
}

# Add two new paragraphs after paragraph 3 ("Antrag: {Antrag}")
$p3 = $d.Paragraphs(3).Range
$p3.InsertParagraphAfter()
Write-Host "After InsertParagraphAfter (1):"
foreach ($p in $d.Paragraphs) {
    Write-Host "PARA: [" $p.Range.Text "]"
}

$d.Paragraphs(4).Range.InsertBefore("-------------------------")
Write-Host "After filling paragraph 4:"
foreach ($p in $d.Paragraphs) {
    Write-Host "PARA: [" $p.Range.Text "]"
}

$d.Paragraphs(4).Range.InsertParagraphAfter()
Write-Host "After InsertParagraphAfter (2):"
foreach ($p in $d.Paragraphs) {
    Write-Host "PARA: [" $p.Range.Text "]"
}

$d.Paragraphs(5).Range.InsertBefore("{/rows}")
Write-Host "Final:"
foreach ($p in $d.Paragraphs) {
    Write-Host "PARA: [" $p.Range.Text "]"
}
